$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 702.9091
$ws.Range("I111").Value = 716.125
$ws.Range("J111").Value = 667.6667
$ws.Range("K111").Value = 2148.375
$ws.Range("L111").Value = 2003.0001
$ws.Range("M111").Value = 918.625
$ws.Range("N111").Value = -8137.0001
$ws.Range("H113").Value = 1667.4286
$ws.Range("I113").Value = 1545
$ws.Range("J113").Value = 1742.7693
$ws.Range("K113").Value = 1545
$ws.Range("L113").Value = 1742.7693
$ws.Range("M113").Value = 1709
$ws.Range("N113").Value = -8250.7693
$ws.Range("H116").Value = 4000
$ws.Range("I116").Value = 2500
$ws.Range("J116").Value = 5800
$ws.Range("K116").Value = 2500
$ws.Range("L116").Value = 5800
$ws.Range("M116").Value = 942
$ws.Range("N116").Value = -12684
$ws.Range("H125").Value = 2563
$ws.Range("I125").Value = 3157.75
$ws.Range("J125").Value = 2087.2
$ws.Range("K125").Value = 28419.75
$ws.Range("L125").Value = 18784.8
$ws.Range("M125").Value = -25959.75
$ws.Range("N125").Value = -23704.8
$ws.Range("H132").Value = 246600.66
$ws.Range("I132").Value = 297143.6
$ws.Range("J132").Value = 1106.4286
$ws.Range("K132").Value = 891430.7999999999
$ws.Range("L132").Value = 3319.2858
$ws.Range("M132").Value = -888900.7999999999
$ws.Range("N132").Value = -8379.2858
$ws.Range("H137").Value = 43479950
$ws.Range("I137").Value = 1628.9474
$ws.Range("J137").Value = 250001970
$ws.Range("K137").Value = 4886.8422
$ws.Range("L137").Value = 750005910
$ws.Range("M137").Value = -2336.8422
$ws.Range("N137").Value = -750011010
$ws.Range("H141").Value = 640.1429000000001
$ws.Range("I141").Value = 601.11536
$ws.Range("J141").Value = 1147.5
$ws.Range("K141").Value = 1803.34608
$ws.Range("L141").Value = 3442.5
$ws.Range("M141").Value = 3376.65392
$ws.Range("N141").Value = -13802.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5123.9243
$ws.Range("I32").Value = 5547.75
$ws.Range("J32").Value = 3549.7144
$ws.Range("K32").Value = 5547.75
$ws.Range("L32").Value = 3549.7144
$ws.Range("M32").Value = -5260.75
$ws.Range("N32").Value = -4123.7144
$ws.Range("H61").Value = 2303.0435
$ws.Range("I61").Value = 1602.8
$ws.Range("J61").Value = 6971.3335
$ws.Range("K61").Value = 1602.8
$ws.Range("L61").Value = 6971.3335
$ws.Range("M61").Value = -1390.8
$ws.Range("N61").Value = -7395.3335
$ws.Range("H129").Value = 49949.5
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 49949.5
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 49949.5
$ws.Range("N129").Value = -59949.5
$ws.Range("H132").Value = 1686.359
$ws.Range("I132").Value = 1639.3793
$ws.Range("J132").Value = 1822.6
$ws.Range("K132").Value = 4918.1379
$ws.Range("L132").Value = 5467.799999999999
$ws.Range("M132").Value = -2388.1379
$ws.Range("N132").Value = -10527.8
$ws.Range("H136").Value = 2303.0435
$ws.Range("I136").Value = 1602.8
$ws.Range("J136").Value = 6971.3335
$ws.Range("K136").Value = 4808.4
$ws.Range("L136").Value = 20914.0005
$ws.Range("M136").Value = -2258.4
$ws.Range("N136").Value = -26014.0005

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H98").Value = 59393
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 59393
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 59393
$ws.Range("N98").Value = -65383
$ws.Range("H99").Value = 960.875
$ws.Range("I99").Value = 715
$ws.Range("J99").Value = 1452.625
$ws.Range("K99").Value = 715
$ws.Range("L99").Value = 1452.625
$ws.Range("M99").Value = 783
$ws.Range("N99").Value = -4448.625
$ws.Range("H109").Value = 59166.668
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 59166.668
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 59166.668
$ws.Range("N109").Value = -61940.668
$ws.Range("H134").Value = 59916.293
$ws.Range("I134").Value = 67751.8
$ws.Range("J134").Value = 1150
$ws.Range("K134").Value = 203255.4
$ws.Range("L134").Value = 3450
$ws.Range("M134").Value = -200720.4
$ws.Range("N134").Value = -8520

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1864.6666
$ws.Range("I31").Value = 1922.75
$ws.Range("J31").Value = 1400
$ws.Range("K31").Value = 1922.75
$ws.Range("L31").Value = 1400
$ws.Range("M31").Value = -1627.75
$ws.Range("N31").Value = -1990
$ws.Range("H34").Value = 1864.6666
$ws.Range("I34").Value = 1922.75
$ws.Range("J34").Value = 1400
$ws.Range("K34").Value = 1922.75
$ws.Range("L34").Value = 1400
$ws.Range("M34").Value = -1720.75
$ws.Range("N34").Value = -1804
$ws.Range("H86").Value = 22731214
$ws.Range("I86").Value = 4777.625
$ws.Range("J86").Value = 35717750
$ws.Range("K86").Value = 4777.625
$ws.Range("L86").Value = 35717750
$ws.Range("M86").Value = -3654.625
$ws.Range("N86").Value = -35719996
$ws.Range("H89").Value = 22731214
$ws.Range("I89").Value = 4777.625
$ws.Range("J89").Value = 35717750
$ws.Range("K89").Value = 23888.125
$ws.Range("L89").Value = 178588750
$ws.Range("M89").Value = -18272.125
$ws.Range("N89").Value = -178599982
$ws.Range("H99").Value = 1387.2307
$ws.Range("I99").Value = 1350.7142
$ws.Range("J99").Value = 1429.8334
$ws.Range("K99").Value = 1350.7142
$ws.Range("L99").Value = 1429.8334
$ws.Range("M99").Value = 147.2858000000001
$ws.Range("N99").Value = -4425.8334
$ws.Range("H126").Value = 1387.2307
$ws.Range("I126").Value = 1350.7142
$ws.Range("J126").Value = 1429.8334
$ws.Range("K126").Value = 4052.1426
$ws.Range("L126").Value = 4289.5002
$ws.Range("M126").Value = -1582.1426
$ws.Range("N126").Value = -9229.5002
$ws.Range("H132").Value = 4184.9
$ws.Range("I132").Value = 3550.4285
$ws.Range("J132").Value = 5665.3335
$ws.Range("K132").Value = 10651.2855
$ws.Range("L132").Value = 16996.0005
$ws.Range("M132").Value = -8121.2855
$ws.Range("N132").Value = -22056.0005
$ws.Range("H134").Value = 2548.8696
$ws.Range("I134").Value = 2689.647
$ws.Range("J134").Value = 2150
$ws.Range("K134").Value = 8068.941
$ws.Range("L134").Value = 6450
$ws.Range("M134").Value = -5533.941
$ws.Range("N134").Value = -11520

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 58824860
$ws.Range("I132").Value = 90910080
$ws.Range("J132").Value = 1950
$ws.Range("K132").Value = 818190720
$ws.Range("L132").Value = 17550
$ws.Range("M132").Value = -818188190
$ws.Range("N132").Value = -22610
$ws.Range("H134").Value = 3488.3684
$ws.Range("I134").Value = 3075.2354
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 9225.706200000001
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -4155.706200000001
$ws.Range("N134").Value = -31140

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3005.1177
$ws.Range("I126").Value = 2949.4167
$ws.Range("J126").Value = 3138.8
$ws.Range("K126").Value = 8848.250100000001
$ws.Range("L126").Value = 9416.400000000001
$ws.Range("M126").Value = -6378.250100000001
$ws.Range("N126").Value = -14356.4
$ws.Range("H132").Value = 2266.6
$ws.Range("I132").Value = 1732.4445
$ws.Range("J132").Value = 3067.8333
$ws.Range("K132").Value = 5197.333500000001
$ws.Range("L132").Value = 9203.499899999999
$ws.Range("M132").Value = -2667.333500000001
$ws.Range("N132").Value = -14263.4999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 887.9818
$ws.Range("I136").Value = 751.3095
$ws.Range("J136").Value = 1329.5385
$ws.Range("K136").Value = 2253.9285
$ws.Range("L136").Value = 3988.6155
$ws.Range("M136").Value = 296.0715
$ws.Range("N136").Value = -9088.6155

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 982.7759
$ws.Range("I132").Value = 446.88373
$ws.Range("J132").Value = 2519
$ws.Range("K132").Value = 1340.65119
$ws.Range("L132").Value = 7557
$ws.Range("M132").Value = 1189.34881
$ws.Range("N132").Value = -12617
